$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

# The card number is a long, purely numeric string. Plain assignment would
# let Excel auto-convert it to a numeric value, but the source workbook
# stores it as text, so force text formatting before writing it.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 10.02.2025"

$ws.Range("B6").Value = "11.02."
$ws.Range("C6").Value = "12.02."
$ws.Range("D6").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E6").Value = "54,52-"

$ws.Range("B7").Value = "14.02."
$ws.Range("C7").Value = "15.02."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 96572705"
$ws.Range("E7").Value = "37,52-"

$ws.Range("B8").Value = "17.02."
$ws.Range("C8").Value = "18.02."
$ws.Range("D8").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E8").Value = "63,47-"

$ws.Range("B9").Value = "19.02."
$ws.Range("C9").Value = "20.02."
$ws.Range("D9").Value = "PAYPAL KHHNCC"
$ws.Range("E9").Value = "89,14-"

$ws.Range("B10").Value = "23.02."
$ws.Range("C10").Value = "24.02."
$ws.Range("D10").Value = "PAYPAL BBPCLS"
$ws.Range("E10").Value = "8,16-"

$ws.Range("B11").Value = "24.02."
$ws.Range("C11").Value = "25.02."
$ws.Range("D11").Value = "BEITRAG Allianz SE K-59919972"
$ws.Range("E11").Value = "52,99-"

$ws.Range("D12").Value = "KONTOSTAND AM 27.02.2025"
$ws.Range("E12").Value = "305,80-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 05.03.2025"
